$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A216").Value = "IMX-USD"
$ws.Range("A217").Value = "TAO-USD"
$ws.Range("A218").Value = "GRT-USD"
